# appoint1: adjusting and updating the appointment form
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("survey")

# Widen column G (calculation column) to fit new content.
# (32.67 "ColumnWidth" round-trips to the OOXML width of 33.5 chars.)
$ws.Columns.Item(7).ColumnWidth = 32.67

# Row 23: "Type" question label becomes "Appointment Type".
$ws.Range("C23").Value = "Appointment Type"

# Insert a new row before row 24 for the new "notes" question,
# pushing lab_test/date_appoint/date_reminder rows down by one.
$ws.Rows.Item(24).Insert()

# Copy the style used by the other survey-row cells (e.g. row 23) onto
# the newly inserted row so it matches the rest of the question rows.
$ws.Range("A23:C23").Copy()
$ws.Range("A24:C24").PasteSpecial(-4122)
$excel.CutCopyMode = 0

$ws.Range("A24").Value = "string"
$ws.Range("B24").Value = "notes"
$ws.Range("C24").Value = "Any notes about this Appointment?"

# The former "date_reminder" row (now shifted to row 27) is removed entirely.
$ws.Rows.Item(27).Delete()
